# Apply updated cryptocurrency price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.804.83"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "'4.002.36"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'531.57"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").Value = "'151.00"
$ws.Range("E6").Value = "  +1.88%  "
$ws.Range("D7").Value = "'0.692"
$ws.Range("E7").Value = "  +10.14%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.744"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  -3.57%  "
$ws.Range("D11").Value = "'0.0000327"
$ws.Range("E11").Value = "  -4.84%  "
$ws.Range("D12").Value = "'47.40"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").Value = "'10.64"
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("D14").Value = "'4.648.00"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").Value = "'4.005.02"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("E16").Value = "  -2.95%  "
$ws.Range("E17").Value = "  -4.41%  "
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").Value = "'71.730.40"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "'425.48"
$ws.Range("E21").Value = "  -4.40%  "
$ws.Range("D22").Value = "'97.38"
$ws.Range("E22").Value = "  +2.62%  "
$ws.Range("E23").Value = "  -3.94%  "
$ws.Range("D24").Value = "'4.20"
$ws.Range("E24").Value = "  +3.31%  "
$ws.Range("D25").Value = "'14.34"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").Value = "'11.25"
$ws.Range("E26").Value = "  -8.87%  "
$ws.Range("D27").Value = "'10.68"
$ws.Range("E27").Value = "  -3.63%  "
$ws.Range("D28").Value = "'5.83"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").Value = "'36.57"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("D30").Value = "'3.58"
$ws.Range("E30").Value = "  +23.06%  "
$ws.Range("D31").Value = "'13.35"
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'677.90"
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'7.03"
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'44.14"
$ws.Range("E35").Value = "  +6.28%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'65.74"
$ws.Range("E36").Value = "  -3.58%  "
$ws.Range("D37").Value = "'0.433"
$ws.Range("E37").Value = "  -4.26%  "
$ws.Range("D38").Value = "'0.152"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").Value = "'0.0₃0823"
$ws.Range("E39").Value = "  -9.01%  "
$ws.Range("D40").Value = "'3.42"
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0486"
$ws.Range("E43").Value = "  -1.99%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'3.22"
$ws.Range("E44").Value = "  +2.98%  "
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("D46").Value = "'3.43"
$ws.Range("E46").Value = "  -3.35%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.63"
$ws.Range("E47").Value = "  -8.32%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "'9.58"
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("E49").Value = "  -6.38%  "
$ws.Range("D50").Value = "'0.000272"
$ws.Range("E50").Value = "  -3.46%  "
$ws.Range("D51").Value = "'145.58"
$ws.Range("E51").Value = "  +1.36%  "
